$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Control Logic for ROM")

# Fill in the ROM control-signal columns (H:O) for the rows that were left
# blank in the starter file. Values are entered as text to match the
# existing "binary digit string" formatting used throughout the table
# (RegWEn, ImmSel, BrUn, ASel, BSel, ALUSel, MemRW, WBSel).

$rows = @{
    17 = @("1","000","0","1","0","0000","0","00")   # lb
    18 = @("1","000","0","1","0","0000","0","00")   # lh
    19 = @("1","000","0","1","0","0000","0","00")   # lw
    28 = @("0","001","0","1","0","0000","1","00")   # sb
    29 = @("0","001","0","1","0","0000","1","00")   # sh
    30 = @("0","001","0","1","0","0000","1","00")   # sw
    37 = @("1","011","0","0","0","0000","0","01")   # auipc
    38 = @("1","011","0","0","0","1111","0","01")   # lui
    39 = @("1","100","0","0","0","0000","0","10")   # jal
    40 = @("1","000","0","1","0","0000","0","10")   # jalr
}

$cols = @("H","I","J","K","L","M","N","O")

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $vals[$i]
    }
}

$wb.Save()
